# The workbook has one worksheet per backward-elimination step, each with a
# large "OLS Regression Results" summary blob (statsmodels console output)
# stored in cell B2. Every one of those blobs contains a line like:
#
#   Time:                        20:51:43   Log-Likelihood:                -174.64
#
# The commit simply re-ran/re-printed the summary a few minutes later, so the
# wall-clock "Time:" stamp moved from 20:51:43 / 20:51:44 to 20:59:45 on every
# sheet, while all other statistics stayed identical. Replace that timestamp
# wherever it appears, leaving the rest of each summary untouched.

$wb = $excel.ActiveWorkbook

$oldTimes = @("20:51:43", "20:51:44")
$newTime = "20:59:45"

foreach ($ws in $wb.Worksheets) {
    $cell = $ws.Range("B2")
    $text = $cell.Value2

    if ($text -eq $null) {
        continue
    }

    if ($text.IndexOf("Time:") -lt 0) {
        continue
    }

    $updated = $text
    foreach ($oldTime in $oldTimes) {
        $updated = $updated.Replace("Time:                        " + $oldTime + "   Log-Likelihood:", "Time:                        " + $newTime + "   Log-Likelihood:")
    }

    if ($updated -ne $text) {
        $cell.Value = $updated
    }
}

Write-Host "Updated Time: stamp to $newTime on all matching sheets"
